# "tested out spacy and decided no"
# Merge the spell-checked / grammar-checked bold heading runs (dropping the
# w:proofErr spellStart/spellEnd/gramStart/gramEnd markers Word's proofer had
# inserted around "Qdrant", "ConversationalRetrievalChain" and "large") back
# into single plain runs, and append the little word-split test-run summary
# at the bottom of the doc.

$d = $word.ActiveDocument

function Merge-BoldHeading {
    param(
        [string]$matchText,
        [string]$finalText
    )
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $ptext = $p.Range.Text.TrimEnd("`r")
        if ($ptext -eq $matchText) {
            # Insert a brand-new (proofErr-free) paragraph ahead of the
            # proofed one, give it the merged text + bold formatting, then
            # delete the old paragraph outright -- that drops every
            # w:proofErr sibling that lived inside it along with it.
            $p.Range.InsertParagraphBefore()
            $newp = $d.Paragraphs.Item($i)
            $newp.Range.Text = $finalText
            $newp.Range.Font.Bold = $true

            $old = $d.Paragraphs.Item($i + 1)
            $old.Range.Delete()
            return
        }
    }
    throw "paragraph not found: $matchText"
}

Merge-BoldHeading "Qdrant tutorial video:" "Qdrant tutorial video:"
Merge-BoldHeading "Qdrant binary-quantization documentation:" "Qdrant binary-quantization documentation:"
Merge-BoldHeading "ConversationalRetrievalChain documentation:" "ConversationalRetrievalChain documentation:"
Merge-BoldHeading "Qdrant documentation to access text-embeddings-3-large" "Qdrant documentation to access text-embeddings-3-large"

# Replace the very last (empty) paragraph with the new test-run write-up,
# one paragraph per line.
$texts = @(
    "Below is test documents processed:",
    "12-page document of job descriptions",
    "2-page resume",
    "18 words split incorrectly",
    "4347 words total",
    ".41%",
    "See: Test docs with word splits"
)

$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = $texts[0]
for ($i = 1; $i -lt $texts.Length; $i++) {
    $p.Range.InsertParagraphAfter()
    $p = $d.Paragraphs.Item($p.Index + 1)
    $p.Range.Text = $texts[$i]
}

Write-Output "done"
